# Commit: "how tf is adx test passing and tr update alone failing HUH?"
#
# The ADX (Average Directional Index) helper columns F:P were only
# populated with real formulas starting a few rows down (once the
# rolling windows had enough history). The early "warm-up" rows were
# missing those cells entirely. This backfills them with literal 0s so
# every row in the used range (A1:P200) has a value in every column,
# and restores the selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: none of F:P exist yet - fill the whole span with zeros.
$ws.Range("F1:P1").Value = 0

# Rows 2-14: F already sits between E (MAX/ABS tr formula) and G
# (directional-movement formulas), and K:P (the smoothed/ADX columns)
# don't exist yet either - both gaps get zero-filled per row.
for ($r = 2; $r -le 14; $r++) {
    $row = $ws.Range("F${r},K${r}:P${r}")
    foreach ($area in $row.Areas) {
        $area.Value = 0
    }
}

# Rows 15-27: only the final ADX column P is missing.
for ($r = 15; $r -le 27; $r++) {
    $ws.Range("P$r").Value = 0
}

# Restore the saved selection/view described by the workbook.
[void]$ws.Range("E2").Select()
